# Generate Report for Handoff
# Refresh the localization status report with a new run:
#  - old file id  869fce92-48a3-4e65-8728-ddba8a4ac588
#  - new file id  d3627a0d-745b-4ce1-9ae3-92e25b28d728
#  - new handoff xliff hash 638d068fd090a30f31bc3a8c579211c0518d986c
#  - refreshed timestamps

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

$oldId = "869fce92-48a3-4e65-8728-ddba8a4ac588"
$newId = "d3627a0d-745b-4ce1-9ae3-92e25b28d728"

$newHandoffZhCn = "$newId.638d068fd090a30f31bc3a8c579211c0518d986c.zh-cn.xlf"
$newHandoffDeDe = "$newId.638d068fd090a30f31bc3a8c579211c0518d986c.de-de.xlf"

# Hyperlink targets stay pinned to the historical commit/file that is
# actually linked to; only the visible display text and the cell values
# refresh to the new run's identifiers.
$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8ed2b79e5ea792ae75d32d0521bc078b79e7609/e2e/$oldId.md"

### ---- Overview sheet ----
$ws_overview.Range("A2").Value = "$newId.md"
$ws_overview.Range("B2").Value = "e2e\$newId.md"
$ws_overview.Range("G2").Value = "2016-08-31 17:08:30"

$ws_overview.Hyperlinks.Delete()
$ws_overview.Hyperlinks.Add($ws_overview.Range("B2"), $hyperlinkTarget, [System.Type]::Missing, [System.Type]::Missing, "e2e\$newId.md")

### ---- zh-cn sheet ----
$ws_zhcn.Range("A2").Value = "$newId.md"
$ws_zhcn.Range("G2").Value = $newHandoffZhCn
$ws_zhcn.Range("H2").Value = "2016-08-31 17:08:26"

$ws_zhcn.Hyperlinks.Delete()
$ws_zhcn.Hyperlinks.Add($ws_zhcn.Range("A2"), $hyperlinkTarget, [System.Type]::Missing, [System.Type]::Missing, "$newId.md")

### ---- de-de sheet ----
$ws_dede.Range("A2").Value = "$newId.md"
$ws_dede.Range("G2").Value = $newHandoffDeDe
$ws_dede.Range("H2").Value = "2016-08-31 17:08:30"

$ws_dede.Hyperlinks.Delete()
$ws_dede.Hyperlinks.Add($ws_dede.Range("A2"), $hyperlinkTarget, [System.Type]::Missing, [System.Type]::Missing, "$newId.md")
